$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: SV-230223 / passes / 1d / passed / 3000-01-01 / John Doe ---
$ws.Range("A2").Value = "SV-230223"
$ws.Range("B2").Value = "This control passes according to this attestation"
$ws.Range("C2").Value = "1d"
$ws.Range("D2").Value = "passed"
$ws.Range("E2").Value = (Get-Date -Year 3000 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("F2").Value = "John Doe"

# --- Row 3: SV-230223 / passes / 1d / passed / 1999-01-01 / John Doe ---
$ws.Range("A3").Value = "SV-230223"
$ws.Range("B3").Value = "This control passes according to this attestation"
$ws.Range("C3").Value = "1d"
$ws.Range("D3").Value = "passed"
$ws.Range("E3").Value = (Get-Date -Year 1999 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("F3").Value = "John Doe"

# --- Row 4: SV-230223 / fails / 1d / failed / 3000-01-01 / John Doe ---
$ws.Range("A4").Value = "SV-230223"
$ws.Range("B4").Value = "This control fails according to this attestation"
$ws.Range("C4").Value = "1d"
$ws.Range("D4").Value = "failed"
$ws.Range("E4").Value = (Get-Date -Year 3000 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("F4").Value = "John Doe"

# --- Row 5: SV-230223 / fails / 1d / failed / 1999-01-01 / John Doe ---
$ws.Range("A5").Value = "SV-230223"
$ws.Range("B5").Value = "This control fails according to this attestation"
$ws.Range("C5").Value = "1d"
$ws.Range("D5").Value = "failed"
$ws.Range("E5").Value = (Get-Date -Year 1999 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("F5").Value = "John Doe"

# --- Row 6: SV-230221 / passes / 1d / passed / 3000-01-01 / John Doe ---
$ws.Range("A6").Value = "SV-230221"
$ws.Range("B6").Value = "This control passes according to this attestation"
$ws.Range("C6").Value = "1d"
$ws.Range("D6").Value = "passed"
$ws.Range("E6").Value = (Get-Date -Year 3000 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("F6").Value = "John Doe"

# --- Row 7: SV-230221 / passes / 1d / passed / 1999-01-01 / John Doe ---
$ws.Range("A7").Value = "SV-230221"
$ws.Range("B7").Value = "This control passes according to this attestation"
$ws.Range("C7").Value = "1d"
$ws.Range("D7").Value = "passed"
$ws.Range("E7").Value = (Get-Date -Year 1999 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("F7").Value = "John Doe"

# --- Row 8: SV-230222 / passes / 1d / passed / 3000-01-01 / John Doe ---
$ws.Range("A8").Value = "SV-230222"
$ws.Range("B8").Value = "This control passes according to this attestation"
$ws.Range("C8").Value = "1d"
$ws.Range("D8").Value = "passed"
$ws.Range("E8").Value = (Get-Date -Year 3000 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("F8").Value = "John Doe"

# --- Row 9: SV-230222 / passes / 1d / passed / 1999-01-01 / John Doe ---
$ws.Range("A9").Value = "SV-230222"
$ws.Range("B9").Value = "This control passes according to this attestation"
$ws.Range("C9").Value = "1d"
$ws.Range("D9").Value = "passed"
$ws.Range("E9").Value = (Get-Date -Year 1999 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("F9").Value = "John Doe"

# --- Row 10: V-73166 / passes / 153d / passed / literal timestamp text / John Doe ---
$ws.Range("A10").Value = "V-73166"
$ws.Range("B10").Value = "This control passes according to this attestation"
$ws.Range("C10").Value = "153d"
$ws.Range("D10").Value = "passed"
# Column E is styled as a date column; this row stores a literal ISO-8601
# timestamp string instead of a date serial, so reset the style before
# writing the text so it does not inherit the date number format.
$ws.Range("E10").Style = "Normal"
$ws.Range("E10").Value = "2024-03-21T22:17:52.761Z"
$ws.Range("F10").Value = "John Doe"

# Match the saved selection/active cell from the edited workbook.
$ws.Range("A4").Select()
